$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC6_SearchResults_Typeahead")
$ws2 = $wb.Worksheets.Item("Testdata")

# --- Sheet1: TC6_SearchResults_Typeahead ---
# Remove the two "CLICK_JS" certificate-handling steps (and their related
# WAIT rows) that are no longer part of the test case, rows 3-7. The rows
# below shift up to close the gap.
$ws1.Rows("3:7").Delete()

# --- Sheet2: Testdata ---
# Leave the data as-is; just restore/normalize the remembered selection.
$ws2.Activate()
$ws2.Range("B6").Select()

# Switch back to the main sheet - it becomes the active/selected tab.
$ws1.Activate()
$ws1.Range("B5").Select()
